$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force Text format, assign, then restore the default "Normal" style so the
# cell keeps the same (unstyled) look as before the edit.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "26.194.73"
$ws.Range("E2").Value = "  -2.04%  "

$ws.Range("D3").Value = "1.582.17"
$ws.Range("E3").Value = "  -1.27%  "

$ws.Range("E4").Value = "  -0.33%  "

Set-TextValue $ws.Range("D5") "209.71"
$ws.Range("E5").Value = "  -1.00%  "

Set-TextValue $ws.Range("D6") "0.497"
$ws.Range("E6").Value = "  -3.10%  "

$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D8") "0.0610"
$ws.Range("E8").Value = "  -1.56%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D9") "0.246"
$ws.Range("E9").Value = "  -0.88%  "

Set-TextValue $ws.Range("D10") "19.52"
$ws.Range("E10").Value = "  -1.25%  "

Set-TextValue $ws.Range("D11") "0.0845"
$ws.Range("E11").Value = "  -0.14%  "

$ws.Range("D12").Value = "1.805.36"
$ws.Range("E12").Value = "  -1.17%  "

$ws.Range("D13").Value = "1.580.26"
$ws.Range("E13").Value = "  -1.37%  "

Set-TextValue $ws.Range("D14") "4.05"
$ws.Range("E14").Value = "  +0.17%  "

Set-TextValue $ws.Range("D15") "0.514"
$ws.Range("E15").Value = "  -1.59%  "

Set-TextValue $ws.Range("D16") "64.49"
$ws.Range("E16").Value = "  -0.93%  "

$ws.Range("D17").Value = "26.199.07"
$ws.Range("E17").Value = "  -1.88%  "

$ws.Range("D18").Value = "0.0₃0733"
$ws.Range("E18").Value = "  -0.90%  "

Set-TextValue $ws.Range("D19") "7.27"
$ws.Range("E19").Value = "  +1.02%  "

$ws.Range("E20").Value = "  -0.38%  "

Set-TextValue $ws.Range("D21") "207.14"
$ws.Range("E21").Value = "  -1.54%  "

$ws.Range("E22").Value = "  -0.67%  "

$ws.Range("E23").Value = "  -3.45%  "

Set-TextValue $ws.Range("D24") "8.89"
$ws.Range("E24").Value = "  -1.07%  "

Set-TextValue $ws.Range("D25") "144.61"
$ws.Range("E25").Value = "  +0.66%  "

$ws.Range("E26").Value = "  -0.46%  "

Set-TextValue $ws.Range("D27") "7.00"
$ws.Range("E27").Value = "  -1.45%  "

Set-TextValue $ws.Range("D28") "0.112"
$ws.Range("E28").Value = "  -1.50%  "

Set-TextValue $ws.Range("D29") "15.22"
$ws.Range("E29").Value = "  -1.18%  "

$ws.Range("E30").Value = "  -1.67%  "

$ws.Range("E31").Value = "  -0.96%  "

Set-TextValue $ws.Range("D32") "3.21"
$ws.Range("E32").Value = "  -2.07%  "

$ws.Range("E33").Value = "  -1.01%  "

$ws.Range("D34").Value = "1.276.30"
$ws.Range("E34").Value = "  -1.50%  "

$ws.Range("E35").Value = "  -0.34%  "

Set-TextValue $ws.Range("D36") "0.612"
$ws.Range("E36").Value = "  +1.53%  "

Set-TextValue $ws.Range("D37") "1.48"
$ws.Range("E37").Value = "  -1.04%  "

$ws.Range("E38").Value = "  -2.01%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D39") "0.816"
$ws.Range("E39").Value = "  -1.86%  "

$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D40") "1.02"
$ws.Range("E40").Value = "  -11.93%  "

$ws.Range("E41").Value = "  +2.51%  "

$ws.Range("E42").Value = "  -2.38%  "

$ws.Range("E43").Value = "  -3.12%  "

Set-TextValue $ws.Range("D44") "62.28"
$ws.Range("E44").Value = "  -1.30%  "

$ws.Range("D45").Value = "1.718.73"

Set-TextValue $ws.Range("D46") "89.08"
$ws.Range("E46").Value = "  -1.66%  "

Set-TextValue $ws.Range("D47") "1.55"
$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("E48").Value = "  -0.95%  "

Set-TextValue $ws.Range("D49") "0.0506"
$ws.Range("E49").Value = "  -2.08%  "

$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("E51").Value = "  +0.74%  "
